$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update crypto price/volume table cells per the latest data refresh.
# A few Price-column values end in a trailing zero after the decimal
# point (e.g. "6.50"); prefix those with a literal leading apostrophe
# so Excel stores them as text instead of silently renormalising them
# to numbers ("6.50" -> 6.5) and dropping the zero.

$ws.Range('D2').Value = '44.135.20'
$ws.Range('E2').Value = '  +2.86%  '
$ws.Range('D3').Value = '2.260.84'
$ws.Range('E3').Value = '  +2.27%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '318.69'
$ws.Range('E5').Value = '  +0.05%  '
$ws.Range('D6').Value = '101.82'
$ws.Range('E6').Value = '  +4.51%  '
$ws.Range('D7').Value = '0.581'
$ws.Range('E7').Value = '  +0.66%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('E9').Value = '  -1.42%  '
$ws.Range('D10').Value = '37.24'
$ws.Range('E10').Value = '  +2.54%  '
$ws.Range('D11').Value = '0.0839'
$ws.Range('E11').Value = '  +2.20%  '
$ws.Range('D12').Value = '7.63'
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('E13').Value = '  -0.12%  '
$ws.Range('D14').Value = '2.606.39'
$ws.Range('E14').Value = '  +2.04%  '
$ws.Range('D15').Value = '14.62'
$ws.Range('E15').Value = '  +2.98%  '
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('D17').Value = '2.266.74'
$ws.Range('E17').Value = '  +2.48%  '
$ws.Range('D18').Value = '44.047.26'
$ws.Range('E18').Value = '  +2.92%  '
$ws.Range('D19').Value = '13.46'
$ws.Range('E19').Value = '  -1.30%  '
$ws.Range('D20').Value = '0.0₃0990'
$ws.Range('E20').Value = '  +3.87%  '
$ws.Range('D21').Value = '''6.50'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('D22').Value = '65.79'
$ws.Range('E22').Value = '  +1.52%  '
$ws.Range('D23').Value = '3.12'
$ws.Range('E23').Value = '  -2.90%  '
$ws.Range('D24').Value = '''235.70'
$ws.Range('E24').Value = '  +0.59%  '
$ws.Range('D25').Value = '2.09'
$ws.Range('E25').Value = '  -3.56%  '
$ws.Range('E26').Value = '  +0.15%  '
$ws.Range('D27').Value = '10.59'
$ws.Range('E27').Value = '  +6.70%  '
$ws.Range('D28').Value = '38.74'
$ws.Range('E28').Value = '  +7.68%  '
$ws.Range('E29').Value = '  -1.36%  '
$ws.Range('E30').Value = '  -1.20%  '
$ws.Range('D31').Value = '162.47'
$ws.Range('E31').Value = '  +4.97%  '
$ws.Range('D32').Value = '20.25'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('E33').Value = '  -0.52%  '
$ws.Range('D34').Value = '2.68'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('E35').Value = '  +8.31%  '
$ws.Range('E36').Value = '  +9.97%  '
$ws.Range('D37').Value = '3.07'
$ws.Range('E37').Value = '  -4.98%  '
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = '16.63'
$ws.Range('E39').Value = '  +20.34%  '
$ws.Range('D40').Value = '''3.70'
$ws.Range('E40').Value = '  +1.54%  '
$ws.Range('D41').Value = '4.21'
$ws.Range('E41').Value = '  -2.88%  '
$ws.Range('E42').Value = '  +0.25%  '
$ws.Range('E43').Value = '  +0.06%  '
$ws.Range('D44').Value = '1.789.35'
$ws.Range('E44').Value = '  +3.66%  '
$ws.Range('D45').Value = '0.198'
$ws.Range('E45').Value = '  -0.89%  '
$ws.Range('B46').Value = 'THORChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D46').Value = '5.23'
$ws.Range('E46').Value = '  +0.59%  '
$ws.Range('B47').Value = 'BitcoinSV'
$ws.Range('C47').Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D47').Value = '81.96'
$ws.Range('E47').Value = '  -1.47%  '
$ws.Range('B48').Value = 'ordi'
$ws.Range('C48').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D48').Value = '74.77'
$ws.Range('E48').Value = '  +3.90%  '
$ws.Range('D49').Value = '104.95'
$ws.Range('E49').Value = '  +3.04%  '
$ws.Range('B50').Value = 'Stacks'
$ws.Range('C50').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D50').Value = '1.68'
$ws.Range('E50').Value = '  +7.77%  '
$ws.Range('B51').Value = 'MultiversX'
$ws.Range('C51').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D51').Value = '58.32'
$ws.Range('E51').Value = '  +2.43%  '

Write-Output "Updated $($ws.Name) with latest cryptos data."
